$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 7
$ws.Range("H7").Value = 2168.3333
$ws.Range("I7").Value = 2168.3333
$ws.Range("K7").Value = 2168.3333
$ws.Range("M7").Value = -2056.3333
# Row 11
$ws.Range("H11").Value = 347.91666
$ws.Range("I11").Value = 347.91666
$ws.Range("K11").Value = 347.91666
$ws.Range("M11").Value = -207.91666
# Row 14
$ws.Range("H14").Value = 2168.3333
$ws.Range("I14").Value = 2168.3333
$ws.Range("K14").Value = 2168.3333
$ws.Range("M14").Value = -1977.3333
# Row 40
$ws.Range("H40").Value = 1998.3334
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 1995
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 1995
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2345
# Row 64
$ws.Range("H64").Value = 6333
$ws.Range("J64").Value = 6333
$ws.Range("L64").Value = 6333
$ws.Range("N64").Value = -6829
# Row 67
$ws.Range("H67").Value = 6333
$ws.Range("J67").Value = 6333
$ws.Range("L67").Value = 6333
$ws.Range("N67").Value = -8049
# Row 100
$ws.Range("H100").Value = 3661.3076
$ws.Range("I100").Value = 3800
$ws.Range("J100").Value = 3649.75
$ws.Range("K100").Value = 3800
$ws.Range("L100").Value = 3649.75
$ws.Range("M100").Value = -3259
$ws.Range("N100").Value = -4731.75
# Row 113
$ws.Range("H113").Value = 3616.3333
$ws.Range("J113").Value = 850
$ws.Range("L113").Value = 850
$ws.Range("N113").Value = -7358
# Row 116
$ws.Range("H116").Value = 4532
$ws.Range("I116").Value = 4399
$ws.Range("J116").Value = 4598.5
$ws.Range("K116").Value = 4399
$ws.Range("L116").Value = 4598.5
$ws.Range("M116").Value = -957
$ws.Range("N116").Value = -11482.5
# Row 131
$ws.Range("H131").Value = 797.6667
$ws.Range("I131").Value = 797.6667
$ws.Range("K131").Value = 2393.0001
$ws.Range("M131").Value = 2646.9999
# Row 132
$ws.Range("H132").Value = 168719.5
$ws.Range("I132").Value = 168719.5
$ws.Range("K132").Value = 506158.5
$ws.Range("M132").Value = -503628.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2613.7144
$ws.Range("I45").Value = 2613.7144
$ws.Range("K45").Value = 2613.7144
$ws.Range("M45").Value = -2236.7144
# Row 61
$ws.Range("H61").Value = 3519.75
$ws.Range("I61").Value = 3519.75
$ws.Range("K61").Value = 3519.75
$ws.Range("M61").Value = -3307.75
# Row 97
$ws.Range("H97").Value = 47619620
$ws.Range("I97").Value = 55556156
$ws.Range("J97").Value = 400
$ws.Range("K97").Value = 55556156
$ws.Range("L97").Value = 400
$ws.Range("M97").Value = -55555660
$ws.Range("N97").Value = -1392
# Row 102
$ws.Range("H102").Value = 9593018
$ws.Range("I102").Value = 689149.4399999999
$ws.Range("J102").Value = 33336666
$ws.Range("K102").Value = 689149.4399999999
$ws.Range("L102").Value = 33336666
$ws.Range("M102").Value = -687527.4399999999
$ws.Range("N102").Value = -33339910
# Row 122
$ws.Range("H122").Value = 1900
$ws.Range("I122").Value = 1900
$ws.Range("K122").Value = 5700
$ws.Range("M122").Value = -3250
# Row 132
$ws.Range("H132").Value = 1908.8889
$ws.Range("I132").Value = 1897.5
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 5692.5
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -3162.5
$ws.Range("N132").Value = -11060
# Row 136
$ws.Range("H136").Value = 3519.75
$ws.Range("I136").Value = 3519.75
$ws.Range("K136").Value = 10559.25
$ws.Range("M136").Value = -8009.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 74724.2
$ws.Range("I94").Value = 92913.586
$ws.Range("K94").Value = 92913.586
$ws.Range("M94").Value = -92462.586
# Row 99
$ws.Range("H99").Value = 4618.0527
$ws.Range("I99").Value = 4359.1875
$ws.Range("J99").Value = 5998.6665
$ws.Range("K99").Value = 4359.1875
$ws.Range("L99").Value = 5998.6665
$ws.Range("M99").Value = -2861.1875
$ws.Range("N99").Value = -8994.666499999999
# Row 102
$ws.Range("H102").Value = 11518.625
$ws.Range("I102").Value = 12607
$ws.Range("J102").Value = 3900
$ws.Range("K102").Value = 12607
$ws.Range("L102").Value = 3900
$ws.Range("M102").Value = -9362
$ws.Range("N102").Value = -10390
# Row 105
$ws.Range("H105").Value = 2611
$ws.Range("I105").Value = 2712.8333
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 2712.8333
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -965.8332999999998
$ws.Range("N105").Value = -5494

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 587.25
$ws.Range("I16").Value = 516.3333
$ws.Range("K16").Value = 516.3333
$ws.Range("M16").Value = -229.3333
# Row 35
$ws.Range("H35").Value = 3139.6
$ws.Range("I35").Value = 3362
$ws.Range("J35").Value = 2250
$ws.Range("K35").Value = 3362
$ws.Range("L35").Value = 2250
$ws.Range("M35").Value = -3068
$ws.Range("N35").Value = -2838
# Row 88
$ws.Range("H88").Value = 24009.916
$ws.Range("J88").Value = 24009.916
$ws.Range("L88").Value = 24009.916
$ws.Range("N88").Value = -24821.916
# Row 91
$ws.Range("H91").Value = 24009.916
$ws.Range("J91").Value = 24009.916
$ws.Range("L91").Value = 24009.916
$ws.Range("N91").Value = -26817.916
# Row 99
$ws.Range("H99").Value = 2040800
$ws.Range("I99").Value = 1668000
$ws.Range("J99").Value = 2600000
$ws.Range("K99").Value = 1668000
$ws.Range("L99").Value = 2600000
$ws.Range("M99").Value = -1666502
$ws.Range("N99").Value = -2602996
# Row 105
$ws.Range("H105").Value = 968.25
$ws.Range("I105").Value = 624.3333
$ws.Range("K105").Value = 624.3333
$ws.Range("M105").Value = 1122.6667
# Row 113
$ws.Range("H113").Value = 587.25
$ws.Range("I113").Value = 516.3333
$ws.Range("K113").Value = 516.3333
$ws.Range("M113").Value = 1653.6667
# Row 126
$ws.Range("H126").Value = 2040800
$ws.Range("I126").Value = 1668000
$ws.Range("J126").Value = 2600000
$ws.Range("K126").Value = 5004000
$ws.Range("L126").Value = 7800000
$ws.Range("M126").Value = -5001530
$ws.Range("N126").Value = -7804940

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Range("H6").Value = 618.5
$ws.Range("I6").Value = 49
$ws.Range("K6").Value = 147
$ws.Range("M6").Value = -34
# Row 129
$ws.Range("H129").Value = 1322
$ws.Range("I129").Value = 769.5
$ws.Range("J129").Value = 1874.5
$ws.Range("K129").Value = 2308.5
$ws.Range("L129").Value = 5623.5
$ws.Range("M129").Value = 2691.5
$ws.Range("N129").Value = -15623.5
# Row 132
$ws.Range("H132").Value = 951.3333
$ws.Range("I132").Value = 849.5
$ws.Range("K132").Value = 7645.5
$ws.Range("M132").Value = -5115.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2692.2
$ws.Range("I80").Value = 2692.2
$ws.Range("K80").Value = 2692.2
$ws.Range("M80").Value = -1694.2
# Row 83
$ws.Range("H83").Value = 2692.2
$ws.Range("I83").Value = 2692.2
$ws.Range("K83").Value = 13461
$ws.Range("M83").Value = -8469
# Row 97
$ws.Range("H97").Value = 599.7143
$ws.Range("I97").Value = 449.75
$ws.Range("K97").Value = 449.75
$ws.Range("M97").Value = 46.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
# Row 68
$ws.Range("H68").Value = 5217.909
$ws.Range("I68").Value = 5999.5
$ws.Range("J68").Value = 4280
$ws.Range("K68").Value = 5999.5
$ws.Range("L68").Value = 4280
$ws.Range("M68").Value = -5250.5
$ws.Range("N68").Value = -5778
# Row 71
$ws.Range("H71").Value = 5217.909
$ws.Range("I71").Value = 5999.5
$ws.Range("J71").Value = 4280
$ws.Range("K71").Value = 29997.5
$ws.Range("L71").Value = 21400
$ws.Range("M71").Value = -26253.5
$ws.Range("N71").Value = -28888
# Row 122
$ws.Range("H122").Value = 3413.5
$ws.Range("I122").Value = 2950.75
$ws.Range("K122").Value = 8852.25
$ws.Range("M122").Value = -6402.25
# Row 136
$ws.Range("H136").Value = 5584
$ws.Range("J136").Value = 7752
$ws.Range("L136").Value = 23256
$ws.Range("N136").Value = -28356

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1460.9231
$ws.Range("I81").Value = 1387.2222
$ws.Range("J81").Value = 1626.75
$ws.Range("K81").Value = 2774.4444
$ws.Range("L81").Value = 3253.5
$ws.Range("M81").Value = -1713.4444
$ws.Range("N81").Value = -5375.5
# Row 84
$ws.Range("H84").Value = 1460.9231
$ws.Range("I84").Value = 1387.2222
$ws.Range("J84").Value = 1626.75
$ws.Range("K84").Value = 13872.222
$ws.Range("L84").Value = 16267.5
$ws.Range("M84").Value = -8568.222
$ws.Range("N84").Value = -26875.5
# Row 122
$ws.Range("H122").Value = 1918.25
$ws.Range("J122").Value = 2465.6667
$ws.Range("L122").Value = 7397.000100000001
$ws.Range("N122").Value = -12297.0001
# Row 126
$ws.Range("H126").Value = 5645.1816
$ws.Range("I126").Value = 4982
$ws.Range("J126").Value = 6441
$ws.Range("K126").Value = 14946
$ws.Range("L126").Value = 19323
$ws.Range("M126").Value = -12476
$ws.Range("N126").Value = -24263
# Row 132
$ws.Range("H132").Value = 2175.875
$ws.Range("I132").Value = 1486.7142
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 4460.142599999999
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -1930.142599999999
$ws.Range("N132").Value = -26060
# Row 136
$ws.Range("H136").Value = 1540.25
$ws.Range("I136").Value = 1540.25
$ws.Range("K136").Value = 4620.75
$ws.Range("M136").Value = -2070.75

